$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp (20240814-104249 -> 20240815-094616)
$ws.Name = "IClientBalance-20240815-094616-"

# Update the "balance date" column (G) for every data row (2-274) from 45518 (2024-08-14)
# to 45519 (2024-08-15). Doing this as one range write keeps it fast and preserves the
# existing date number format (style index) already applied to those cells.
$ws.Range("G2:G274").Value = 45519

# A handful of rows also got their balance amounts (columns E and H, which mirror each
# other) corrected for the new day.
$ws.Range("E58").Value = 0
$ws.Range("H58").Value = 0

$ws.Range("E102").Value = 28352.18
$ws.Range("H102").Value = 28352.18

$ws.Range("E105").Value = 16874.68
$ws.Range("H105").Value = 16874.68

$ws.Range("E231").Value = 28881.77
$ws.Range("H231").Value = 28881.77

$ws.Range("E255").Value = 16209.88
$ws.Range("H255").Value = 16209.88

# Move the active selection from B3 to B2, matching the saved cursor position in the file.
$ws.Range("B2").Select()
